$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (shifts existing rows 18-48 down to 19-49)
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with the new data record
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = 45175
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = 100112022
$ws.Cells.Item(18, 7).Value = "Arveja Verde"
$ws.Cells.Item(18, 8).Value = "Perfection"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 50
$ws.Cells.Item(18, 11).Value = 27000
$ws.Cells.Item(18, 12).Value = 28000
$ws.Cells.Item(18, 13).Value = 27600
$ws.Cells.Item(18, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 1104
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
